$wb = $excel.ActiveWorkbook

# --- Sheet "y": append new monthly data points (rows 50-61) ---
$wsY = $wb.Worksheets.Item("y")

$yData = @(
    @(44256, 7001),
    @(44287, 7332),
    @(44317, 8201),
    @(44348, 6932),
    @(44378, 5988),
    @(44409, 6450),
    @(44440, 7923),
    @(44470, 6021),
    @(44501, 5302),
    @(44531, 5542),
    @(44562, 6450),
    @(44593, 7203)
)

$row = 50
foreach ($pair in $yData) {
    $wsY.Range("A$row").Value = $pair[0]
    $wsY.Range("B$row").Value = $pair[1]
    $row = $row + 1
}

# --- Sheet "X": append new monthly data points (rows 62-73) ---
$wsX = $wb.Worksheets.Item("X")

$xData = @(
    @(44621, 3000, 6893),
    @(44652, 3121, 6751),
    @(44682, 3240, 6678),
    @(44713, 2425, 6816),
    @(44743, 2555, 6724),
    @(44774, 2345, 6677),
    @(44805, 2745, 6344),
    @(44835, 2464, 6433),
    @(44866, 2334, 6441),
    @(44896, 2131, 6454),
    @(44927, 1998, 6321),
    @(44958, 2034, 6212)
)

$row = 62
foreach ($triple in $xData) {
    $wsX.Range("A$row").Value = $triple[0]
    $wsX.Range("B$row").Value = $triple[1]
    $wsX.Range("C$row").Value = $triple[2]
    $row = $row + 1
}

# --- Update the selection left behind on each sheet ---
[void]$wsY.Range("F47").Select()
[void]$wsX.Range("G36").Select()

$wsInfo = $wb.Worksheets.Item("Info")
$wsInfo.Activate()
[void]$wsInfo.Range("I15").Select()

# --- Drop the internal "CorpoS" confidentiality header stamp on every sheet ---
foreach ($name in @("Info", "y", "X")) {
    $sheet = $wb.Worksheets.Item($name)
    $ps = $sheet.PageSetup
    $ps.LeftHeader = ""
    $ps.CenterHeader = ""
    $ps.RightHeader = ""
}
